$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.905.46"
$ws.Range("E2").Value = "  +0.85%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.108.23"
$ws.Range("E3").Value = "  +2.48%  "
$ws.Range("E4").Value = "  +0.57%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "559.12"
$ws.Range("E5").Value = "  +3.14%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.53"
$ws.Range("E6").Value = "  +1.73%  "
$ws.Range("E7").Value = "  +0.21%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.098.63"
$ws.Range("E8").Value = "  +2.21%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.489"
$ws.Range("E9").Value = "  +1.76%  "
$ws.Range("E10").Value = "  +5.70%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.151"
$ws.Range("E11").Value = "  +1.64%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.456"
$ws.Range("E12").Value = "  +2.10%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "35.63"
$ws.Range("E13").Value = "  +2.22%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000217"
$ws.Range("E14").Value = "  +1.74%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.610.19"
$ws.Range("E15").Value = "  +3.08%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.017.75"
$ws.Range("E16").Value = "  +1.35%  "
$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.111"
$ws.Range("E17").Value = "  +0.79%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.118.24"
$ws.Range("E18").Value = "  +3.49%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "505.33"
$ws.Range("E19").Value = "  +5.93%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.68"
$ws.Range("E20").Value = "  +3.05%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.73"
$ws.Range("E21").Value = "  +2.72%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.704"
$ws.Range("E22").Value = "  +4.66%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.25"
$ws.Range("E23").Value = "  +2.87%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.44"
$ws.Range("E24").Value = "  +2.58%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "77.95"
$ws.Range("E25").Value = "  +1.21%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.79"
$ws.Range("E27").Value = "  +5.51%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.47"
$ws.Range("E28").Value = "  +8.75%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.06"
$ws.Range("E29").Value = "  +1.06%  "
$ws.Range("E30").Value = "  +0.53%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "26.23"
$ws.Range("E31").Value = "  +3.63%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.55"
$ws.Range("E32").Value = "  -0.02%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.12"
$ws.Range("E33").Value = "  +0.82%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "547.00"
$ws.Range("E34").Value = "  -3.86%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "55.53"
$ws.Range("E35").Value = "  +7.06%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.92"
$ws.Range("E36").Value = "  +1.28%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.24"
$ws.Range("E37").Value = "  -1.05%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0417"
$ws.Range("E38").Value = "  +6.53%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0803"
$ws.Range("E39").Value = "  +3.72%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.066.36"
$ws.Range("E40").Value = "  +5.54%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.119"
$ws.Range("E41").Value = "  +3.19%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.12"
$ws.Range("E42").Value = "  +0.71%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.62"
$ws.Range("E43").Value = "  -6.28%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.257"
$ws.Range("E44").Value = "  +8.24%  "
$ws.Range("B45").Value = "USDe"
$ws.Range("C45").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("E45").Value = "  +0.05%  "
$ws.Range("B46").Value = "Fetch.AI"
$ws.Range("C46").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.12"
$ws.Range("E46").Value = "  +3.72%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "120.74"
$ws.Range("E47").Value = "  +2.24%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "24.45"
$ws.Range("E48").Value = "  +0.92%  "
$ws.Range("E49").Value = "  +0.38%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0₃0504"
$ws.Range("E50").Value = "  -1.96%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.04"
$ws.Range("E51").Value = "  +1.23%  "
